# Apply "Error Calculations and Plots" edit:
# - Update various data cells (some values added, some cleared) across rows 3-33.
# - Row 26 ("RM 232") is removed and all subsequent rows shift up by one,
#   with a new "SC 232" row ending up at row 33.
# - Former rows 34 ("SC 193") and 35 ("SC 232") are dropped, shrinking the
#   sheet's used range from A1:F35 down to A1:F33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..33, columns A..F (ID, A, B, C, D, F)
$data = @(
    @("RM 2", -19.7, 14.9, -13.5, -7.2, 18.03),
    @("RM 8", -19.7, $null, -14.2, -5.7, 17.64),
    @("RM 9", -18.7, 11, -15.4, -6.4, 17.97),
    @("RM 14", -19.5, 12.3, -14.4, $null, 17.66),
    @("RM 21", -19.8, 15.1, -14.2, -5.7, 16.43),
    @("RM 32", -19.5, 15, -13.8, -7.1, 17.24),
    @("RM 38", -19.9, 15.5, -13.9, -6.6, 17.05),
    @("RM 42", $null, 10.5, -14.5, -6.8, 17.26),
    @("RM 52 a", -19.8, 11.5, -14.7, -6.1, 16.43),
    @("RM 58", -20.8, 11.4, -15.5, -7.9, 17.65),
    @("RM 81", -19.9, 12.5, -14.1, $null, $null),
    @("RM 88", -19.9, 12.5, -13.9, -5.3, 17.1),
    @("RM 90", -19.6, 14.4, -13.1, $null, 17.76),
    @("RM 95", -19.1, 12.5, -15.2, -8.4, 16.2),
    @("RM 103", -19.5, 13.5, -15.3, -6.9, 17.34),
    @("RM 116", -19.4, 11.2, -14.7, -7.3, $null),
    @("RM 120", $null, 11.5, -15.2, $null, 18.35),
    @("RM 125", -20.6, 13.2, -15.5, $null, $null),
    @("RM 134", $null, 12.5, -14, -7.2, 17.73),
    @("RM 135", -18.9, 12.7, -14.3, -8.7, 16.58),
    @("RM 138", -19.3, 13.5, -15.4, -6.1, 16.81),
    @("RM 140", -19.5, 12.2, -13.9, -7, 16.48),
    @("RM 142a", -17.7, 12.7, -13.9, -8.1, 16.78),
    @("RM 145", -19.5, 10.7, -15.5, -7.1, 16.6),
    @("SC 5", -20.2, 10.8, -13.8, -5, 17.38),
    @("SC 101", $null, 10, -14.6, -10, 17),
    @("SC 105", -19.6, 11.1, -13.7, -5.9, 17.44),
    @("SC 119", -19.5, 11.2, -13, $null, 18.06),
    @("SC 120", -19.7, 11.4, -13.6, -5.7, 16.89),
    @("SC 132", -18.8, 15.3, -13.7, -8.1, $null),
    @("SC 193", -19.9, 10.5, -14.7, -6.4, $null),
    @("SC 232", -19.5, 10.4, -14.1, -10.7, 17.53)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Remove the now-obsolete trailing rows (old rows 34 and 35), which shrinks
# the sheet from A1:F35 to A1:F33.
$lastRow = $startRow + $data.Count - 1
$ws.Range($ws.Cells.Item($lastRow + 1, 1), $ws.Cells.Item(35, 6)).Clear()
